$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Make room for the new rows -----------------------------------------
# Current data occupies rows 2-7. The refreshed scrape inserts one brand new
# row at the very top (row 2), keeps the former row 2 where it lands (row 3),
# then inserts two more brand new rows (rows 4-5) before the remaining old
# rows continue (old rows 3-7 end up at rows 6-10).
$ws.Range("A2").EntireRow.Insert()
$ws.Range("A4:A5").EntireRow.Insert()

# --- 2. Write the brand new postings ----------------------------------------
$ws.Range("A2").Value = "2025-12-05 18:24:53"
$ws.Range("B2").Value = "【IT導入補助金A類型】受注管理DX「パンダリンク」開発+店舗LP制作+顧客LINE構築の依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5448148"
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = "◆開発 ◇管理"

$ws.Range("A4").Value = "2025-12-05 18:24:53"
$ws.Range("B4").Value = "【急募】不動産向けWordPressとLINE統合の専門家を探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5448323"
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = "○WordPress"

$ws.Range("A5").Value = "2025-12-05 18:24:53"
$ws.Range("B5").Value = "【急募】不動産向けWordPressとLINE統合の専門家を探しています!"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5448342"
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = "○WordPress"

# --- 3. Refresh the "seen at" timestamp for every row still in the list ----
$ws.Range("A3").Value = "2025-12-05 18:24:53"
$ws.Range("A6").Value = "2025-12-05 18:24:53"
$ws.Range("A7").Value = "2025-12-05 18:24:53"
$ws.Range("A8").Value = "2025-12-05 18:24:53"
$ws.Range("A9").Value = "2025-12-05 18:24:53"
$ws.Range("A10").Value = "2025-12-05 18:24:53"

# --- 4. Rebuild the URL hyperlinks ------------------------------------------
# Row inserts do not relocate hyperlink relationships in this object model,
# so drop every existing hyperlink and recreate them all at their final
# (post-insert) positions.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5448148")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5437832")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5448323")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5448342")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5447772")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5447970")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5341051")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5447771")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5447991")

# Re-apply the workbook's existing "Hyperlink" cell style so every F-column
# cell keeps using the same shared style that was already in the file.
$ws.Range("F2:F10").Style = "Hyperlink"
